# "Atualizando o arquivo XLSX" — refresh of the FlashScore odds sheet for
# the 2024-10-28 round. Re-applies the updated odds for several existing
# matches, drops the Peru "Comerciantes Unidos - AD Tarma" fixture, and
# adds the new Colombia "Ind. Medellin - Jaguares de Cordoba" fixture in
# its place (shifting the untouched Paraguay fixture from row 10 to row 11).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Gimnasia L.P. - Union de Santa Fe (odds refresh) ---
$ws.Range("O2").Value = 1.62
$ws.Range("P2").Value = 2.2
$ws.Range("Q2").Value = 3.1
$ws.Range("R2").Value = 1.36

# --- Row 4: Oriente Petrolero - SA Bulo Bulo (odds refresh) ---
$ws.Range("Q4").Value = 1.53
$ws.Range("R4").Value = 2.4

# --- Row 5: Cuiaba - Corinthians (odds refresh) ---
$ws.Range("O5").Value = 1.57
$ws.Range("P5").Value = 2.38

# --- Row 7: Coritiba - CRB (odds refresh) ---
$ws.Range("O7").Value = 1.36
$ws.Range("P7").Value = 3
$ws.Range("Q7").Value = 2.2
$ws.Range("R7").Value = 1.65

# --- Row 13: FC Cincinnati - New York City (odds refresh) ---
$ws.Range("H13").Value = 3.7
$ws.Range("J13").Value = 2.5
$ws.Range("K13").Value = 2.38
$ws.Range("S13").Value = 1.3
$ws.Range("T13").Value = 3.4
$ws.Range("U13").Value = 1.57
$ws.Range("V13").Value = 2.25
$ws.Range("W13").Value = 10
$ws.Range("Y13").Value = 8.5
$ws.Range("AA13").Value = 15
$ws.Range("AC13").Value = 15
$ws.Range("AF13").Value = 41
$ws.Range("AM13").Value = 126
$ws.Range("AO13").Value = 10
$ws.Range("AS13").Value = 101
$ws.Range("AT13").Value = 3.4
$ws.Range("AY13").Value = 23
$ws.Range("BC13").Value = 401

# --- Row 14: Seattle Sounders - Houston Dynamo (odds refresh) ---
$ws.Range("G14").Value = 1.85
$ws.Range("H14").Value = 3.5
$ws.Range("I14").Value = 4.5
$ws.Range("J14").Value = 2.5
$ws.Range("L14").Value = 4.75
$ws.Range("AH14").Value = 21
$ws.Range("AO14").Value = 10
$ws.Range("AX14").Value = 23
$ws.Range("BA14").Value = 101
$ws.Range("BC14").Value = 126

# --- Row 10: new fixture - COLOMBIA PRIMERA A, Ind. Medellin v Jaguares de Cordoba ---
# (replaces the Paraguay fixture that used to live here; it now moves to row 11)
$ws.Range("A10").Value = "K4qymNjt"
$ws.Range("C10").Value = "22:15"
$ws.Range("D10").Value = "COLOMBIA - PRIMERA A"
$ws.Range("E10").Value = "Ind. Medellin"
$ws.Range("F10").Value = "Jaguares de Cordoba"
$ws.Range("G10").Value = 1.45
$ws.Range("H10").Value = 3.9
$ws.Range("I10").Value = 8
$ws.Range("J10").Value = 2.05
$ws.Range("K10").Value = 2.2
$ws.Range("L10").Value = 7.5
$ws.Range("M10").Value = 1.07
$ws.Range("N10").Value = 9
$ws.Range("O10").Value = 1.36
$ws.Range("P10").Value = 3
$ws.Range("Q10").Value = 2.15
$ws.Range("R10").Value = 1.67
$ws.Range("S10").Value = 1.44
$ws.Range("T10").Value = 2.63
$ws.Range("U10").Value = 2.38
$ws.Range("V10").Value = 1.53
$ws.Range("W10").Value = 5.5
$ws.Range("X10").Value = 6
$ws.Range("Y10").Value = 9
$ws.Range("Z10").Value = 9.5
$ws.Range("AA10").Value = 15
$ws.Range("AB10").Value = 34
$ws.Range("AC10").Value = 8
$ws.Range("AD10").Value = 8
$ws.Range("AE10").Value = 23
$ws.Range("AF10").Value = 81
$ws.Range("AG10").Value = 15
$ws.Range("AH10").Value = 41
$ws.Range("AI10").Value = 23
$ws.Range("AJ10").Value = 101
$ws.Range("AK10").Value = 67
$ws.Range("AL10").Value = 67
$ws.Range("AM10").Value = 201
$ws.Range("AN10").Value = 3.2
$ws.Range("AO10").Value = 7.5
$ws.Range("AP10").Value = 23
$ws.Range("AQ10").Value = 23
$ws.Range("AR10").Value = 51
$ws.Range("AS10").Value = 201
$ws.Range("AT10").Value = 2.63
$ws.Range("AU10").Value = 10
$ws.Range("AV10").Value = 81
$ws.Range("AW10").Value = 8.5
$ws.Range("AX10").Value = 41
$ws.Range("AY10").Value = 41
$ws.Range("AZ10").Value = 201
$ws.Range("BA10").Value = 251
$ws.Range("BB10").Value = 501
$ws.Range("BC10").Value = 126
$ws.Range("BD10").Value = 126

# --- Row 11: PARAGUAY PRIMERA DIVISION, Guarani v Sp. Luqueno ---
# (this fixture used to be row 10; it moves down to row 11 with a handful
# of its odds refreshed, replacing the removed Peru fixture)
$ws.Range("A11").Value = "0SbJHVr2"
$ws.Range("C11").Value = "19:30"
$ws.Range("D11").Value = "PARAGUAY - PRIMERA DIVISION"
$ws.Range("E11").Value = "Guarani"
$ws.Range("F11").Value = "Sp. Luqueno"
$ws.Range("G11").Value = 2.25
$ws.Range("H11").Value = 3
$ws.Range("I11").Value = 3.4
$ws.Range("J11").Value = 2.5
$ws.Range("K11").Value = 2.1
$ws.Range("L11").Value = 5
$ws.Range("M11").Value = 1.08
$ws.Range("N11").Value = 8
$ws.Range("O11").Value = 1.4
$ws.Range("P11").Value = 2.75
$ws.Range("Q11").Value = 2.35
$ws.Range("R11").Value = 1.57
$ws.Range("S11").Value = 1.44
$ws.Range("T11").Value = 2.63
$ws.Range("U11").Value = 2
$ws.Range("V11").Value = 1.73
$ws.Range("W11").Value = 6.5
$ws.Range("X11").Value = 10
$ws.Range("Y11").Value = 9.5
$ws.Range("Z11").Value = 21
$ws.Range("AA11").Value = 21
$ws.Range("AB11").Value = 34
$ws.Range("AC11").Value = 7
$ws.Range("AD11").Value = 6
$ws.Range("AE11").Value = 17
$ws.Range("AF11").Value = 51
$ws.Range("AG11").Value = 8.5
$ws.Range("AH11").Value = 15
$ws.Range("AI11").Value = 13
$ws.Range("AJ11").Value = 41
$ws.Range("AK11").Value = 34
$ws.Range("AL11").Value = 41
$ws.Range("AM11").Value = 351
$ws.Range("AN11").Value = 3.75
$ws.Range("AO11").Value = 9.5
$ws.Range("AP11").Value = 21
$ws.Range("AQ11").Value = 34
$ws.Range("AR11").Value = 51
$ws.Range("AS11").Value = 151
$ws.Range("AT11").Value = 2.63
$ws.Range("AU11").Value = 8.5
$ws.Range("AV11").Value = 51
$ws.Range("AW11").Value = 6.5
$ws.Range("AX11").Value = 26
$ws.Range("AY11").Value = 34
$ws.Range("AZ11").Value = 81
$ws.Range("BA11").Value = 101
$ws.Range("BB11").Value = 251
$ws.Range("BC11").Value = 51
$ws.Range("BD11").Value = 51
